$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 163 (shifts existing rows 163.. down by one)
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row with the new command entry
$ws.Range("A163").Value = "TXT2PLINE"
$ws.Range("B163").Value = "A utility to explode TEXT and MTEXT to polylines"

# The row insertion shifts the previously-sorted range (A480:A514 -> A481:A515).
# Re-apply the sort on the shifted range so the worksheet's stored sortState
# metadata reflects the new row numbers (data is already in sorted order).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A481:A515"))
$sortObj.SetRange($ws.Range("A481:A515"))
$sortObj.Header = 2
$sortObj.Apply()

# Leave the selection on the newly-added row, mirroring where the author
# was working when the change was saved.
$null = $ws.Range("B164").Select()
